$wb = $excel.ActiveWorkbook

# mmWave(InBed): append rows 64-75 (Value column E is text "In Bed")
$ws = $wb.Worksheets.Item("mmWave(InBed)")
$ws.Cells.Item(64, 1).NumberFormat = "@"
$ws.Cells.Item(64, 1).Value = "2026-02-01"
$ws.Cells.Item(64, 2).Value = "20:17:29"
$ws.Cells.Item(64, 3).Value = "20:00"
$ws.Cells.Item(64, 4).Value = "Bedroom"
$ws.Cells.Item(64, 5).Value = "In Bed"
$ws.Cells.Item(64, 6).Value = "Occupied"

$ws.Cells.Item(65, 1).NumberFormat = "@"
$ws.Cells.Item(65, 1).Value = "2026-02-01"
$ws.Cells.Item(65, 2).Value = "20:17:31"
$ws.Cells.Item(65, 3).Value = "20:00"
$ws.Cells.Item(65, 4).Value = "Bedroom"
$ws.Cells.Item(65, 5).Value = "In Bed"
$ws.Cells.Item(65, 6).Value = "Occupied"

$ws.Cells.Item(66, 1).NumberFormat = "@"
$ws.Cells.Item(66, 1).Value = "2026-02-01"
$ws.Cells.Item(66, 2).Value = "20:17:32"
$ws.Cells.Item(66, 3).Value = "20:00"
$ws.Cells.Item(66, 4).Value = "Bedroom"
$ws.Cells.Item(66, 5).Value = "In Bed"
$ws.Cells.Item(66, 6).Value = "Occupied"

$ws.Cells.Item(67, 1).NumberFormat = "@"
$ws.Cells.Item(67, 1).Value = "2026-02-01"
$ws.Cells.Item(67, 2).Value = "20:17:34"
$ws.Cells.Item(67, 3).Value = "20:00"
$ws.Cells.Item(67, 4).Value = "Bedroom"
$ws.Cells.Item(67, 5).Value = "In Bed"
$ws.Cells.Item(67, 6).Value = "Occupied"

$ws.Cells.Item(68, 1).NumberFormat = "@"
$ws.Cells.Item(68, 1).Value = "2026-02-01"
$ws.Cells.Item(68, 2).Value = "20:17:35"
$ws.Cells.Item(68, 3).Value = "20:00"
$ws.Cells.Item(68, 4).Value = "Bedroom"
$ws.Cells.Item(68, 5).Value = "In Bed"
$ws.Cells.Item(68, 6).Value = "Occupied"

$ws.Cells.Item(69, 1).NumberFormat = "@"
$ws.Cells.Item(69, 1).Value = "2026-02-01"
$ws.Cells.Item(69, 2).Value = "20:17:40"
$ws.Cells.Item(69, 3).Value = "20:00"
$ws.Cells.Item(69, 4).Value = "Bedroom"
$ws.Cells.Item(69, 5).Value = "In Bed"
$ws.Cells.Item(69, 6).Value = "Occupied"

$ws.Cells.Item(70, 1).NumberFormat = "@"
$ws.Cells.Item(70, 1).Value = "2026-02-01"
$ws.Cells.Item(70, 2).Value = "20:17:42"
$ws.Cells.Item(70, 3).Value = "20:00"
$ws.Cells.Item(70, 4).Value = "Bedroom"
$ws.Cells.Item(70, 5).Value = "In Bed"
$ws.Cells.Item(70, 6).Value = "Occupied"

$ws.Cells.Item(71, 1).NumberFormat = "@"
$ws.Cells.Item(71, 1).Value = "2026-02-01"
$ws.Cells.Item(71, 2).Value = "20:17:45"
$ws.Cells.Item(71, 3).Value = "20:00"
$ws.Cells.Item(71, 4).Value = "Bedroom"
$ws.Cells.Item(71, 5).Value = "In Bed"
$ws.Cells.Item(71, 6).Value = "Occupied"

$ws.Cells.Item(72, 1).NumberFormat = "@"
$ws.Cells.Item(72, 1).Value = "2026-02-01"
$ws.Cells.Item(72, 2).Value = "20:18:20"
$ws.Cells.Item(72, 3).Value = "20:00"
$ws.Cells.Item(72, 4).Value = "Bedroom"
$ws.Cells.Item(72, 5).Value = "In Bed"
$ws.Cells.Item(72, 6).Value = "Occupied"

$ws.Cells.Item(73, 1).NumberFormat = "@"
$ws.Cells.Item(73, 1).Value = "2026-02-01"
$ws.Cells.Item(73, 2).Value = "20:18:22"
$ws.Cells.Item(73, 3).Value = "20:00"
$ws.Cells.Item(73, 4).Value = "Bedroom"
$ws.Cells.Item(73, 5).Value = "In Bed"
$ws.Cells.Item(73, 6).Value = "Occupied"

$ws.Cells.Item(74, 1).NumberFormat = "@"
$ws.Cells.Item(74, 1).Value = "2026-02-01"
$ws.Cells.Item(74, 2).Value = "20:18:25"
$ws.Cells.Item(74, 3).Value = "20:00"
$ws.Cells.Item(74, 4).Value = "Bedroom"
$ws.Cells.Item(74, 5).Value = "In Bed"
$ws.Cells.Item(74, 6).Value = "Occupied"

$ws.Cells.Item(75, 1).NumberFormat = "@"
$ws.Cells.Item(75, 1).Value = "2026-02-01"
$ws.Cells.Item(75, 2).Value = "20:18:27"
$ws.Cells.Item(75, 3).Value = "20:00"
$ws.Cells.Item(75, 4).Value = "Bedroom"
$ws.Cells.Item(75, 5).Value = "In Bed"
$ws.Cells.Item(75, 6).Value = "Occupied"

# mmWave(BR): append rows 60-71 (Value column E is numeric)
$ws = $wb.Worksheets.Item("mmWave(BR)")
$ws.Cells.Item(60, 1).NumberFormat = "@"
$ws.Cells.Item(60, 1).Value = "2026-02-01"
$ws.Cells.Item(60, 2).Value = "20:17:30"
$ws.Cells.Item(60, 3).Value = "20:00"
$ws.Cells.Item(60, 4).Value = "Bedroom"
$ws.Cells.Item(60, 5).Value = 5
$ws.Cells.Item(60, 6).Value = "Occupied"

$ws.Cells.Item(61, 1).NumberFormat = "@"
$ws.Cells.Item(61, 1).Value = "2026-02-01"
$ws.Cells.Item(61, 2).Value = "20:17:31"
$ws.Cells.Item(61, 3).Value = "20:00"
$ws.Cells.Item(61, 4).Value = "Bedroom"
$ws.Cells.Item(61, 5).Value = 2
$ws.Cells.Item(61, 6).Value = "Occupied"

$ws.Cells.Item(62, 1).NumberFormat = "@"
$ws.Cells.Item(62, 1).Value = "2026-02-01"
$ws.Cells.Item(62, 2).Value = "20:17:33"
$ws.Cells.Item(62, 3).Value = "20:00"
$ws.Cells.Item(62, 4).Value = "Bedroom"
$ws.Cells.Item(62, 5).Value = 20
$ws.Cells.Item(62, 6).Value = "Occupied"

$ws.Cells.Item(63, 1).NumberFormat = "@"
$ws.Cells.Item(63, 1).Value = "2026-02-01"
$ws.Cells.Item(63, 2).Value = "20:17:35"
$ws.Cells.Item(63, 3).Value = "20:00"
$ws.Cells.Item(63, 4).Value = "Bedroom"
$ws.Cells.Item(63, 5).Value = 2
$ws.Cells.Item(63, 6).Value = "Occupied"

$ws.Cells.Item(64, 1).NumberFormat = "@"
$ws.Cells.Item(64, 1).Value = "2026-02-01"
$ws.Cells.Item(64, 2).Value = "20:17:36"
$ws.Cells.Item(64, 3).Value = "20:00"
$ws.Cells.Item(64, 4).Value = "Bedroom"
$ws.Cells.Item(64, 5).Value = 1
$ws.Cells.Item(64, 6).Value = "Occupied"

$ws.Cells.Item(65, 1).NumberFormat = "@"
$ws.Cells.Item(65, 1).Value = "2026-02-01"
$ws.Cells.Item(65, 2).Value = "20:17:41"
$ws.Cells.Item(65, 3).Value = "20:00"
$ws.Cells.Item(65, 4).Value = "Bedroom"
$ws.Cells.Item(65, 5).Value = 11
$ws.Cells.Item(65, 6).Value = "Occupied"

$ws.Cells.Item(66, 1).NumberFormat = "@"
$ws.Cells.Item(66, 1).Value = "2026-02-01"
$ws.Cells.Item(66, 2).Value = "20:17:43"
$ws.Cells.Item(66, 3).Value = "20:00"
$ws.Cells.Item(66, 4).Value = "Bedroom"
$ws.Cells.Item(66, 5).Value = 2
$ws.Cells.Item(66, 6).Value = "Occupied"

$ws.Cells.Item(67, 1).NumberFormat = "@"
$ws.Cells.Item(67, 1).Value = "2026-02-01"
$ws.Cells.Item(67, 2).Value = "20:17:46"
$ws.Cells.Item(67, 3).Value = "20:00"
$ws.Cells.Item(67, 4).Value = "Bedroom"
$ws.Cells.Item(67, 5).Value = 1
$ws.Cells.Item(67, 6).Value = "Occupied"

$ws.Cells.Item(68, 1).NumberFormat = "@"
$ws.Cells.Item(68, 1).Value = "2026-02-01"
$ws.Cells.Item(68, 2).Value = "20:18:21"
$ws.Cells.Item(68, 3).Value = "20:00"
$ws.Cells.Item(68, 4).Value = "Bedroom"
$ws.Cells.Item(68, 5).Value = 4
$ws.Cells.Item(68, 6).Value = "Occupied"

$ws.Cells.Item(69, 1).NumberFormat = "@"
$ws.Cells.Item(69, 1).Value = "2026-02-01"
$ws.Cells.Item(69, 2).Value = "20:18:23"
$ws.Cells.Item(69, 3).Value = "20:00"
$ws.Cells.Item(69, 4).Value = "Bedroom"
$ws.Cells.Item(69, 5).Value = 2
$ws.Cells.Item(69, 6).Value = "Occupied"

$ws.Cells.Item(70, 1).NumberFormat = "@"
$ws.Cells.Item(70, 1).Value = "2026-02-01"
$ws.Cells.Item(70, 2).Value = "20:18:26"
$ws.Cells.Item(70, 3).Value = "20:00"
$ws.Cells.Item(70, 4).Value = "Bedroom"
$ws.Cells.Item(70, 5).Value = 3
$ws.Cells.Item(70, 6).Value = "Occupied"

$ws.Cells.Item(71, 1).NumberFormat = "@"
$ws.Cells.Item(71, 1).Value = "2026-02-01"
$ws.Cells.Item(71, 2).Value = "20:18:28"
$ws.Cells.Item(71, 3).Value = "20:00"
$ws.Cells.Item(71, 4).Value = "Bedroom"
$ws.Cells.Item(71, 5).Value = 2
$ws.Cells.Item(71, 6).Value = "Occupied"

# mmWave(HR): append rows 60-71 (Value column E is numeric)
$ws = $wb.Worksheets.Item("mmWave(HR)")
$ws.Cells.Item(60, 1).NumberFormat = "@"
$ws.Cells.Item(60, 1).Value = "2026-02-01"
$ws.Cells.Item(60, 2).Value = "20:17:30"
$ws.Cells.Item(60, 3).Value = "20:00"
$ws.Cells.Item(60, 4).Value = "Bedroom"
$ws.Cells.Item(60, 5).Value = 53
$ws.Cells.Item(60, 6).Value = "Occupied"

$ws.Cells.Item(61, 1).NumberFormat = "@"
$ws.Cells.Item(61, 1).Value = "2026-02-01"
$ws.Cells.Item(61, 2).Value = "20:17:31"
$ws.Cells.Item(61, 3).Value = "20:00"
$ws.Cells.Item(61, 4).Value = "Bedroom"
$ws.Cells.Item(61, 5).Value = 50
$ws.Cells.Item(61, 6).Value = "Occupied"

$ws.Cells.Item(62, 1).NumberFormat = "@"
$ws.Cells.Item(62, 1).Value = "2026-02-01"
$ws.Cells.Item(62, 2).Value = "20:17:33"
$ws.Cells.Item(62, 3).Value = "20:00"
$ws.Cells.Item(62, 4).Value = "Bedroom"
$ws.Cells.Item(62, 5).Value = 68
$ws.Cells.Item(62, 6).Value = "Occupied"

$ws.Cells.Item(63, 1).NumberFormat = "@"
$ws.Cells.Item(63, 1).Value = "2026-02-01"
$ws.Cells.Item(63, 2).Value = "20:17:34"
$ws.Cells.Item(63, 3).Value = "20:00"
$ws.Cells.Item(63, 4).Value = "Bedroom"
$ws.Cells.Item(63, 5).Value = 50
$ws.Cells.Item(63, 6).Value = "Occupied"

$ws.Cells.Item(64, 1).NumberFormat = "@"
$ws.Cells.Item(64, 1).Value = "2026-02-01"
$ws.Cells.Item(64, 2).Value = "20:17:36"
$ws.Cells.Item(64, 3).Value = "20:00"
$ws.Cells.Item(64, 4).Value = "Bedroom"
$ws.Cells.Item(64, 5).Value = 49
$ws.Cells.Item(64, 6).Value = "Occupied"

$ws.Cells.Item(65, 1).NumberFormat = "@"
$ws.Cells.Item(65, 1).Value = "2026-02-01"
$ws.Cells.Item(65, 2).Value = "20:17:41"
$ws.Cells.Item(65, 3).Value = "20:00"
$ws.Cells.Item(65, 4).Value = "Bedroom"
$ws.Cells.Item(65, 5).Value = 59
$ws.Cells.Item(65, 6).Value = "Occupied"

$ws.Cells.Item(66, 1).NumberFormat = "@"
$ws.Cells.Item(66, 1).Value = "2026-02-01"
$ws.Cells.Item(66, 2).Value = "20:17:42"
$ws.Cells.Item(66, 3).Value = "20:00"
$ws.Cells.Item(66, 4).Value = "Bedroom"
$ws.Cells.Item(66, 5).Value = 50
$ws.Cells.Item(66, 6).Value = "Occupied"

$ws.Cells.Item(67, 1).NumberFormat = "@"
$ws.Cells.Item(67, 1).Value = "2026-02-01"
$ws.Cells.Item(67, 2).Value = "20:17:45"
$ws.Cells.Item(67, 3).Value = "20:00"
$ws.Cells.Item(67, 4).Value = "Bedroom"
$ws.Cells.Item(67, 5).Value = 49
$ws.Cells.Item(67, 6).Value = "Occupied"

$ws.Cells.Item(68, 1).NumberFormat = "@"
$ws.Cells.Item(68, 1).Value = "2026-02-01"
$ws.Cells.Item(68, 2).Value = "20:18:21"
$ws.Cells.Item(68, 3).Value = "20:00"
$ws.Cells.Item(68, 4).Value = "Bedroom"
$ws.Cells.Item(68, 5).Value = 52
$ws.Cells.Item(68, 6).Value = "Occupied"

$ws.Cells.Item(69, 1).NumberFormat = "@"
$ws.Cells.Item(69, 1).Value = "2026-02-01"
$ws.Cells.Item(69, 2).Value = "20:18:22"
$ws.Cells.Item(69, 3).Value = "20:00"
$ws.Cells.Item(69, 4).Value = "Bedroom"
$ws.Cells.Item(69, 5).Value = 50
$ws.Cells.Item(69, 6).Value = "Occupied"

$ws.Cells.Item(70, 1).NumberFormat = "@"
$ws.Cells.Item(70, 1).Value = "2026-02-01"
$ws.Cells.Item(70, 2).Value = "20:18:26"
$ws.Cells.Item(70, 3).Value = "20:00"
$ws.Cells.Item(70, 4).Value = "Bedroom"
$ws.Cells.Item(70, 5).Value = 51
$ws.Cells.Item(70, 6).Value = "Occupied"

$ws.Cells.Item(71, 1).NumberFormat = "@"
$ws.Cells.Item(71, 1).Value = "2026-02-01"
$ws.Cells.Item(71, 2).Value = "20:18:27"
$ws.Cells.Item(71, 3).Value = "20:00"
$ws.Cells.Item(71, 4).Value = "Bedroom"
$ws.Cells.Item(71, 5).Value = 50
$ws.Cells.Item(71, 6).Value = "Occupied"
